# Project2_Data.xlsx update
# Adds two new "L2 cache" summary tables (Trace1 / Trace2) below the
# existing L1-cache tables, and widens columns C/D to fit the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Column widths: split the old merged C:D "bestFit" width into two
# explicit custom widths.
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 11.666666666666666
$ws.Columns.Item(4).ColumnWidth = 12.166666666666666

# ---------------------------------------------------------------------
# Table 1: Trace1 (rows 74-83)
# ---------------------------------------------------------------------
$ws.Range("B74").Value = "Trace1"
$ws.Range("C74").Value = "512KB, 4-Way"
$ws.Range("D74").Value = "512KB, 8-Way"
$ws.Range("E74").Value = "1MB, 8-Way"

$ws.Range("B75").Value = "L2 Read Accesses"
$ws.Range("C75").Value = 7176477
$ws.Range("D75").Value = 7176477
$ws.Range("E75").Value = 7176477

$ws.Range("B76").Value = "L2 Write Accesses"
$ws.Range("C76").Value = 1399342
$ws.Range("D76").Value = 1399342
$ws.Range("E76").Value = 1399342

$ws.Range("B77").Value = "L2 Cache Accesses"
$ws.Range("C77").Formula = "=C75+C76"
$ws.Range("D77").Formula = "=D75+D76"
$ws.Range("E77").Formula = "=E75+E76"

$ws.Range("B78").Value = "L2 Read Misses"
$ws.Range("C78").Value = 1716334
$ws.Range("D78").Value = 1704921
$ws.Range("E78").Value = 1521667

$ws.Range("B79").Value = "L2 Write Misses"
$ws.Range("C79").Value = 78219
$ws.Range("D79").Value = 82513
$ws.Range("E79").Value = 63213

$ws.Range("B80").Value = "L2 Misses"
$ws.Range("C80").Formula = "=C79+C78"
$ws.Range("D80").Formula = "=D79+D78"
$ws.Range("E80").Formula = "=E79+E78"

$ws.Range("B81").Value = "L2 Miss Rate"
$ws.Range("C81").Formula = "=C80/C77"
$ws.Range("D81").Formula = "=D80/D77"
$ws.Range("E81").Formula = "=E80/E77"

$ws.Range("B83").Value = "Total Execution Time"
$ws.Range("C83").Value = 414757805
$ws.Range("D83").Value = 447965961
$ws.Range("E83").Value = 430900059

# ---------------------------------------------------------------------
# Table 2: Trace2 (rows 85-94)
# ---------------------------------------------------------------------
$ws.Range("B85").Value = "Trace2"
$ws.Range("C85").Value = "512KB, 4-Way"
$ws.Range("D85").Value = "512KB, 8-Way"
$ws.Range("E85").Value = "1MB, 8-Wat"

$ws.Range("B86").Value = "L2 Read Accesses"
$ws.Range("C86").Value = 15108476
$ws.Range("D86").Value = 15108476
$ws.Range("E86").Value = 15108476

$ws.Range("B87").Value = "L2 Write Accesses"
$ws.Range("C87").Value = 5205108
$ws.Range("D87").Value = 5205108
$ws.Range("E87").Value = 5205108

$ws.Range("B88").Value = "L2 Cache Accesses"
$ws.Range("C88").Formula = "=C86+C87"
$ws.Range("D88").Formula = "=D86+D87"
$ws.Range("E88").Formula = "=E86+E87"

$ws.Range("B89").Value = "L2 Read Misses"
$ws.Range("C89").Value = 11063239
$ws.Range("D89").Value = 11174627
$ws.Range("E89").Value = 8721449

$ws.Range("B90").Value = "L2 Write Misses"
$ws.Range("C90").Value = 3927207
$ws.Range("D90").Value = 3935890
$ws.Range("E90").Value = 2921195

$ws.Range("B91").Value = "L2 Misses"
$ws.Range("C91").Formula = "=C90+C89"
$ws.Range("D91").Formula = "=D90+D89"
$ws.Range("E91").Formula = "=E90+E89"

$ws.Range("B92").Value = "L2 Miss Rate"
$ws.Range("C92").Formula = "=C91/C88"
$ws.Range("D92").Formula = "=D91/D88"
$ws.Range("E92").Formula = "=E91/E88"

$ws.Range("B94").Value = "Total Execution Time"
$ws.Range("C94").Value = 1844992985
$ws.Range("D94").Value = 2066107269
$ws.Range("E94").Value = 1612970779

# ---------------------------------------------------------------------
# Recalculate formulas so cached <v> results are written out.
# ---------------------------------------------------------------------
$excel.Calculate()

# ---------------------------------------------------------------------
# Sheet view: scroll down and select I91, matching the final cursor
# position left by the author.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 71
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I91").Select()
